$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "date" and "percentage" looking values in this sheet (From/To columns,
# Success Rate column) are stored as plain text in the workbook, not as
# real Excel dates / percentages. Mark those ranges as Text ("@") before
# writing so Excel's auto-detection doesn't silently convert them into date
# serials / numeric percentages.
$ws.Range("D2:E3").NumberFormat = "@"
$ws.Range("P2:P3").NumberFormat = "@"

# --- Row 2 (Test # 1) updates ---
$ws.Range("C2").Value = "BTCUSDT"
$ws.Range("D2").Value = "2021-09-01"
$ws.Range("E2").Value = "2022-01-01"
$ws.Range("L2").Value = "Early MACD"
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 26
$ws.Range("P2").Value = "42.3%"
$ws.Range("Q2").Value = -4
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 6600
$ws.Range("T2").Value = -6000
$ws.Range("U2").Value = 278.35
$ws.Range("V2").Value = 321.65

# --- Row 3 (Test # 2) updates ---
$ws.Range("B3").Value = "Binance"
$ws.Range("C3").Value = "BTCUSDT"
$ws.Range("D3").Value = "2021-09-01"
$ws.Range("E3").Value = "2022-01-01"
$ws.Range("L3").Value = "Early MACD"
$ws.Range("M3").Value = 10
$ws.Range("N3").Value = 16
$ws.Range("O3").Value = 26
$ws.Range("P3").Value = "38.5%"
$ws.Range("Q3").Value = -6
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 6000
$ws.Range("T3").Value = -6400
$ws.Range("U3").Value = 287.9
$ws.Range("V3").Value = -687.9
